# SCD0011-004 - Penambahan Leads Prospek
# Commit: "Update Excel SCD0011 until SCD0016"
#
# The whole SCD* test-suite workbook set was renumbered (SCD0173 -> SCD0011,
# part of a new SCD0011..SCD0016 range) and the TC_ID on this sheet was
# updated from the old Jira tag "DGS-188" to the new case id "SCD0011-004".
# Column B was also widened (it now holds the longer id string) and the
# cell selection that was active when the file was last saved moved from
# L2 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0173 -> SCD0011
$ws.Name = "SCD0011"

# TC_ID cell: old ticket tag "DGS-188" -> new test case id "SCD0011-004"
$ws.Range("B2").Value = "SCD0011-004"

# Column B was manually resized in Excel to fit the new, longer id text
$ws.Columns.Item(2).ColumnWidth = 12.42578125

# The active cell/selection when the workbook was last saved was B3
$ws.Range("B3").Select()
